# Master_Fuel_Sector_List.xlsx edit
#
# Commit: "Added sector to 3 sector maps, reverted 4 mod-H files to
# pre-CEDS_Data to run Make"
#
# For this workbook the meaningful change is on the "Sectors" sheet: a
# new sector row ("2L1_Oil-tanker-loading") is inserted right before the
# existing "2L_Other-process-emissions" row (row 43), pushing that row
# and everything below it down by one. The new row reuses the same
# "pop" activity / "kt" units / "NC" type pattern already used by its
# neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new row above the current row 43 ("2L_Other-process-emissions"),
# shifting it (and rows 44-59) down to 44-60.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row with the new sector's data.
$ws.Range("A43").Value = "2L1_Oil-tanker-loading"
$ws.Range("B43").Value = "pop"
$ws.Range("C43").Value = "kt"
$ws.Range("D43").Value = "NC"

# Match the saved selection/viewport from the edited workbook.
$ws.Activate() | Out-Null
$ws.Range("B46").Select() | Out-Null
